$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 39; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    $v = $cell.Value2()
    if ($v -eq 46061) {
        $cell.Value = 46062
    }
}
